$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# Fix for enforcing consistent upload spreadsheet value casing
$ws.Range("AA4").Value = "pcr"
$ws.Range("AA5").Value = "Latex AGGLUTINATION, PCR"
$ws.Range("AA7").Value = "lancefield"

# Newly populated example cells
$ws.Range("P4").Value = "m"
$ws.Range("AG6").Value = 1
$ws.Range("AH6").Value = "etest"

# Widen column AA to fit the longer example value
$ws.Columns("AA").ColumnWidth = 51.5

# Restore the frozen-pane scroll position / active selection
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("AD4").Select()
$ws.Range("AI16").Select()
